$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before the existing "District" column (F) to hold
# the school "Address" -- this shifts the old F column (District/Name
# header/Mandya values) one column to the right, into G.
$ws.Columns("F:F").Insert()

# New header for the inserted column.
$ws.Range("F2").Value = "Address"

# For every data row, derive the school address from the second line of
# the "Names" column (column B), which looks like:
#   "<Teacher Name>
#    <School>, <Taluk>, <District>."
# The address is that second line with the trailing ", <District>."
# dropped and the remaining ", " separators removed (matching how the
# Address column values were produced for this sheet).
for ($row = 4; $row -le 33; $row++) {
    $nameCell = $ws.Range("B$row").Value()
    $lines = $nameCell.Split("`n")
    if ($lines.Count -gt 1) {
        $address = $lines[1]
        $address = $address.Replace(", Mandya.", "")
        $address = $address.Replace(", ", "")
        $ws.Range("F$row").Value = $address
    }
}
